$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.175.47'
$ws.Range("E2").Value = '  +0.38%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.422.01'
$ws.Range("E3").Value = '  +0.07%  '

# Row 4
$ws.Range("E4").Value = '  -0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.22'
$ws.Range("E5").Value = '  +1.18%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.84'
$ws.Range("E6").Value = '  +0.54%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.621'
$ws.Range("E7").Value = '  -1.67%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.723'
$ws.Range("E9").Value = '  -0.83%  '

# Row 10
$ws.Range("E10").Value = '  +0.82%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.63'
$ws.Range("E11").Value = '  +0.41%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.29'
$ws.Range("E12").Value = '  +2.86%  '

# Row 13
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.972.57'
$ws.Range("E13").Value = '  +0.43%  '

# Row 14
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000215'
$ws.Range("E14").Value = '  +6.77%  '

# Row 15
$ws.Range("E15").Value = '  -0.53%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.44'
$ws.Range("E16").Value = '  -3.24%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.427.88'
$ws.Range("E17").Value = '  +0.51%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.63'
$ws.Range("E18").Value = '  +2.98%  '

# Row 19
$ws.Range("E19").Value = '  -0.67%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.209.58'
$ws.Range("E20").Value = '  +0.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '466.29'
$ws.Range("E21").Value = '  +4.88%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '90.71'
$ws.Range("E22").Value = '  -0.68%  '

# Row 23
$ws.Range("E23").Value = '  +3.10%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.46'
$ws.Range("E24").Value = '  +4.32%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.44'
$ws.Range("E25").Value = '  +19.39%  '

# Row 26
$ws.Range("E26").Value = '  +2.29%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.02'
$ws.Range("E27").Value = '  -0.06%  '

# Row 28
$ws.Range("E28").Value = '  -0.20%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.64'
$ws.Range("E29").Value = '  +1.27%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.93'
$ws.Range("E30").Value = '  +0.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.64'
$ws.Range("E31").Value = '  -3.20%  '

# Row 32
$ws.Range("E32").Value = '  -1.95%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.111'
$ws.Range("E33").Value = '  -2.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '40.58'
$ws.Range("E34").Value = '  -4.94%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.04%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '59.00'
$ws.Range("E36").Value = '  +10.83%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0486'
$ws.Range("E37").Value = '  -2.14%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.00'
$ws.Range("E38").Value = '  +0.16%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.05'
$ws.Range("E39").Value = '  +4.39%  '

# Row 40
$ws.Range("E40").Value = '  +3.98%  '

# Row 41
$ws.Range("E41").Value = '  +0.10%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.33'
$ws.Range("E42").Value = '  -1.36%  '

# Row 43
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.69'
$ws.Range("E43").Value = '  +11.59%  '

# Row 44
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '145.39'
$ws.Range("E44").Value = '  +3.09%  '

# Row 45
$ws.Range("E45").Value = '  +5.33%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.29'
$ws.Range("E46").Value = '  +1.83%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").Value = '  +18.41%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.39'
$ws.Range("E48").Value = '  -0.71%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.35'
$ws.Range("E49").Value = '  -0.35%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₃0516'
$ws.Range("E50").Value = '  +26.98%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.138'
$ws.Range("E51").Value = '  +4.82%  '
